# Add a "name" column to the materials table on the "mat" sheet, and
# populate it with material names for the four existing materials.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mat")

# Insert a new column before column B ("unit weight"). Everything from the
# old column B onward (option, c, f, c/p, r-elev, piezo, std-dev columns,
# and the helper lookup tables in N:O) shifts one column to the right.
$ws.Columns("B:B").Insert()

# The two header cells that used to read "Mat" / "Piezo" are retyped in
# lower case to match the "mat"/"piezo" sheet-name strings used elsewhere
# in the workbook.
$ws.Range("A3").Value = "mat"
$ws.Range("I3").Value = "piezo"

# New header + data for the inserted "name" column.
$ws.Range("B3").Value = "name"
$ws.Range("B4").Value = "Shell"
$ws.Range("B5").Value = "Core"
$ws.Range("B6").Value = "Clay"
$ws.Range("B7").Value = "Sand"

# Give the new column a wider width so the material names are readable.
$ws.Columns("B:B").ColumnWidth = 16.50203125

# The conditional formatting rules and data-validation lists still point at
# the pre-insert ranges/columns, so repoint them at the shifted locations.
$fc = $ws.Cells.FormatConditions
for ($i = 1; $i -le $fc.Count; $i++) {
    $cond = $fc.Item($i)
    $addr = $cond.AppliesTo.Address()
    $formula = $cond.Formula1
    if ($addr -eq "`$D`$4:`$E`$13") {
        $cond.Formula1 = '=$D4="cp"'
        $ws.Range("E4:F13").FormatConditions.Item($i).AppliesTo = $ws.Range("E4:F13")
    }
}

$ws.Range("D4:E13").FormatConditions.Delete()
$ws.Range("F4:G13").FormatConditions.Delete()
$ws.Range("J4:K13").FormatConditions.Delete()
$ws.Range("L4:L13").FormatConditions.Delete()

$rule1 = $ws.Range("E4:F13").FormatConditions.Add(2, 0, '=$D4="cp"')
$rule2 = $ws.Range("G4:H13").FormatConditions.Add(2, 0, '=$D4="mc"')
$rule3 = $ws.Range("K4:L13").FormatConditions.Add(2, 0, '=$D4="cp"')
$rule4 = $ws.Range("M4:M13").FormatConditions.Add(2, 0, '=$D4="mc"')

# Data validation lists: re-point them at the shifted source ranges.
$ws.Range("I4:I13").Validation.Delete()
$ws.Range("I4:I13").Validation.Add(3, 1, 1, "=`$O`$9:`$O`$10")

$ws.Range("D4:D13").Validation.Delete()
$ws.Range("D4:D13").Validation.Add(3, 1, 1, "=`$O`$5:`$O`$6")

Write-Output "done"
